$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 41 with the new LeetCode entry
$ws.Range("A41").Value = 3440
$ws.Range("B41").Value = "Reschedule Meetings for Maximum Free Time 2"
$ws.Range("C41").Value = "Math"
$ws.Range("D41").Value = "Collect all gaps[], find the max gap to the left/right of the current gap, consider if we can expand the gap or not at the current meeting either to the left or right."

# Update sheet view: new topLeftCell and selection seen in the diff
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("D40").Select()
